# Update 想去人数 (F column) and 最低票价 (G column) figures
# across the four sheets to the refreshed scrape values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 263
$ws.Range("F3").Value = 2809
$ws.Range("G3").Value = 75
$ws.Range("G5").Value = 60
$ws.Range("F7").Value = 3026
$ws.Range("F8").Value = 1921
$ws.Range("F9").Value = 242
$ws.Range("F11").Value = 2582
$ws.Range("F12").Value = 582
$ws.Range("F13").Value = 278
$ws.Range("F14").Value = 14
$ws.Range("F18").Value = 9636
$ws.Range("F22").Value = 7613
$ws.Range("F23").Value = 12156
$ws.Range("F29").Value = 2738
$ws.Range("F32").Value = 2738
$ws.Range("F33").Value = 1191
$ws.Range("F36").Value = 61
$ws.Range("F38").Value = 1165
$ws.Range("F42").Value = 584

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 74
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 88
$ws.Range("F4").Value = 3
$ws.Range("F9").Value = 1191

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 639
$ws.Range("F4").Value = 210

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 639
$ws.Range("F3").Value = 263
$ws.Range("F4").Value = 2809
$ws.Range("G4").Value = 75
$ws.Range("F5").Value = 74
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 88
$ws.Range("G7").Value = 60
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 3026
$ws.Range("F12").Value = 1921
$ws.Range("F14").Value = 242
$ws.Range("F15").Value = 2582
$ws.Range("F17").Value = 582
$ws.Range("F18").Value = 278
$ws.Range("F19").Value = 14
$ws.Range("F23").Value = 9636
$ws.Range("F26").Value = 7613
$ws.Range("F27").Value = 12156
$ws.Range("F33").Value = 2738
$ws.Range("F38").Value = 61
$ws.Range("F46").Value = 584

